$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.99000000000062"
$ws.Range("H2").Value = [double]"6.057010099147497e-06"
$ws.Range("I2").Value = [double]"6.057010099147497e-06"
$ws.Range("L2").Value = [double]"55.75189959573235"
$ws.Range("M2").Value = "[29.738893750290686, 81.76490544117402]"
$ws.Range("N2").Value = [double]"8.587054962827523e-05"
$ws.Range("O2").Value = [double]"8.587054962827523e-05"
$ws.Range("P2").Value = [double]"1.981184556317888"
$ws.Range("Q2").Value = "[1.4654476241970418, 2.496921488438735]"
$ws.Range("R2").Value = [double]"8.341698443103951e-10"
$ws.Range("S2").Value = [double]"8.341698443103951e-10"
$ws.Range("T2").Value = [double]"67.78906701511309"
$ws.Range("U2").Value = "[53.057205724303685, 82.5209283059225]"
$ws.Range("V2").Value = [double]"5.323519403077626e-12"
$ws.Range("W2").Value = [double]"5.323519403077626e-12"
$ws.Range("X2").Value = [double]"17.79495495495538"
$ws.Range("Y2").Value = [double]"15.66164164164202"
$ws.Range("Z2").Value = [double]"19.92826826826874"
$ws.Range("F3").Value = [double]"25.99000000000062"
$ws.Range("H3").Value = [double]"4.263501187451446e-05"
$ws.Range("I3").Value = [double]"4.263501187451446e-05"
$ws.Range("L3").Value = [double]"41.28555549840154"
$ws.Range("M3").Value = "[21.61077542697012, 60.96033556983295]"
$ws.Range("N3").Value = [double]"0.0001144918681363905"
$ws.Range("O3").Value = [double]"0.0001144918681363905"
$ws.Range("P3").Value = [double]"1.427710775505272"
$ws.Range("Q3").Value = "[0.861658045128733, 1.9937635058818106]"
$ws.Range("R3").Value = [double]"7.044243359066726e-06"
$ws.Range("S3").Value = [double]"7.044243359066726e-06"
$ws.Range("T3").Value = [double]"60.10646442365576"
$ws.Range("U3").Value = "[48.06571554349004, 72.14721330382147]"
$ws.Range("V3").Value = [double]"4.389821839367869e-13"
$ws.Range("W3").Value = [double]"4.389821839367869e-13"
$ws.Range("X3").Value = [double]"20.08436436436484"
$ws.Range("Y3").Value = [double]"17.74292292292335"
$ws.Range("Z3").Value = [double]"22.42580580580634"
$ws.Range("F4").Value = [double]"25.99000000000062"
$ws.Range("H4").Value = [double]"5.137887753026149e-07"
$ws.Range("I4").Value = [double]"5.137887753026149e-07"
$ws.Range("L4").Value = [double]"52.93978578011221"
$ws.Range("M4").Value = "[31.098288122016484, 74.78128343820794]"
$ws.Range("N4").Value = [double]"1.362037535823646e-05"
$ws.Range("O4").Value = [double]"1.362037535823646e-05"
$ws.Range("P4").Value = [double]"0.5975001042863468"
$ws.Range("Q4").Value = "[0.16981581911296129, 1.0251843894597323]"
$ws.Range("R4").Value = [double]"0.007234478205516925"
$ws.Range("S4").Value = [double]"0.007234478205516925"
$ws.Range("T4").Value = [double]"56.92643290047059"
$ws.Range("U4").Value = "[45.497298169159265, 68.35556763178192]"
$ws.Range("V4").Value = [double]"4.705125178361413e-13"
$ws.Range("W4").Value = [double]"4.705125178361413e-13"
$ws.Range("X4").Value = [double]"23.51847847847904"
$ws.Range("Y4").Value = [double]"21.74938938938991"
$ws.Range("Z4").Value = [double]"25.28756756756817"
$ws.Range("F5").Value = [double]"25.99000000000062"
$ws.Range("H5").Value = [double]"0.000887690086261661"
$ws.Range("I5").Value = [double]"0.000887690086261661"
$ws.Range("L5").Value = [double]"40.95747589843733"
$ws.Range("M5").Value = "[15.175472862575191, 66.73947893429947]"
$ws.Range("N5").Value = [double]"0.002522823303553245"
$ws.Range("O5").Value = [double]"0.002522823303553245"
$ws.Range("P5").Value = [double]"0.4591316590831935"
$ws.Range("Q5").Value = "[-0.2201316173686534, 1.1383949355350405]"
$ws.Range("R5").Value = [double]"0.1801702110084309"
$ws.Range("S5").Value = [double]"0.1801702110084309"
$ws.Range("T5").Value = [double]"51.57401542087204"
$ws.Range("U5").Value = "[37.785460785304785, 65.3625700564393]"
$ws.Range("V5").Value = [double]"1.660075854559295e-09"
$ws.Range("W5").Value = [double]"1.660075854559295e-09"
$ws.Range("X5").Value = [double]"24.09083083083141"
$ws.Range("Y5").Value = [double]"21.28110110110161"
$ws.Range("Z5").Value = [double]"26.9005605605612"
$ws.Range("F6").Value = [double]"23.54000000000024"
$ws.Range("H6").Value = [double]"3.839671729455318e-07"
$ws.Range("I6").Value = [double]"3.839671729455318e-07"
$ws.Range("L6").Value = [double]"55.6140775628056"
$ws.Range("M6").Value = "[33.66971654990307, 77.55843857570814]"
$ws.Range("N6").Value = [double]"6.493144864228029e-06"
$ws.Range("O6").Value = [double]"6.493144864228029e-06"
$ws.Range("P6").Value = [double]"-0.3522105877898465"
$ws.Range("Q6").Value = "[-0.8050527720910781, 0.10063159651138509]"
$ws.Range("R6").Value = [double]"0.1242317158919686"
$ws.Range("S6").Value = [double]"0.1242317158919686"
$ws.Range("T6").Value = [double]"64.2098381437019"
$ws.Range("U6").Value = "[51.879761908015425, 76.53991437938838]"
$ws.Range("V6").Value = [double]"1.143529715363911e-13"
$ws.Range("W6").Value = [double]"1.143529715363911e-13"
$ws.Range("X6").Value = [double]"1.319559559559572"
$ws.Range("Y6").Value = [double]"-0.3770170170170226"
$ws.Range("Z6").Value = [double]"3.016136136136167"
$ws.Range("F7").Value = [double]"23.54000000000024"
$ws.Range("H7").Value = [double]"0.586556030580129"
$ws.Range("I7").Value = [double]"0.586556030580129"
$ws.Range("L7").Value = [double]"11.08330939697693"
$ws.Range("M7").Value = "[-18.01788171315833, 40.18450050711218]"
$ws.Range("N7").Value = [double]"0.4470401914856685"
$ws.Range("O7").Value = [double]"0.4470401914856685"
$ws.Range("P7").Value = [double]"0.6603948521059619"
$ws.Range("Q7").Value = "[-2.4717635893108896, 3.7925532935228135]"
$ws.Range("R7").Value = [double]"0.6731078855633443"
$ws.Range("S7").Value = [double]"0.6731078855633443"
$ws.Range("T7").Value = [double]"57.53996632521"
$ws.Range("U7").Value = "[42.32121674905164, 72.75871590136836]"
$ws.Range("V7").Value = [double]"1.259604420766891e-09"
$ws.Range("W7").Value = [double]"1.259604420766891e-09"
$ws.Range("X7").Value = [double]"21.06582582582604"
$ws.Range("Y7").Value = [double]"9.331171171171265"
$ws.Range("Z7").Value = [double]"32.80048048048082"
$ws.Range("F8").Value = [double]"23.54000000000024"
$ws.Range("H8").Value = [double]"8.139166135512355e-06"
$ws.Range("I8").Value = [double]"8.139166135512355e-06"
$ws.Range("L8").Value = [double]"54.5511339757145"
$ws.Range("M8").Value = "[32.21079593631603, 76.89147201511297]"
$ws.Range("N8").Value = [double]"1.207779598044745e-05"
$ws.Range("O8").Value = [double]"1.207779598044745e-05"
$ws.Range("P8").Value = [double]"-0.006289474781961957"
$ws.Range("Q8").Value = "[-0.5220264069028095, 0.5094474573388856]"
$ws.Range("R8").Value = [double]"0.9805127398549558"
$ws.Range("S8").Value = [double]"0.9805127398549558"
$ws.Range("T8").Value = [double]"60.57930095554254"
$ws.Range("U8").Value = "[46.44095686562211, 74.71764504546297]"
$ws.Range("V8").Value = [double]"4.25224300215632e-11"
$ws.Range("W8").Value = [double]"4.25224300215632e-11"
$ws.Range("X8").Value = [double]"0.0235635635635667"
$ws.Range("Y8").Value = [double]"-1.908648648648668"
$ws.Range("Z8").Value = [double]"1.955775775775801"

Write-Output "done"
